$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H136").Value = 220000
$ws.Range("J136").Value = 220000
$ws.Range("L136").Value = 220000
$ws.Range("N136").Value = -230200
$ws.Range("H137").Value = 41670060
$ws.Range("I137").Value = 27780628
$ws.Range("J137").Value = 83338350
$ws.Range("K137").Value = 83341884
$ws.Range("L137").Value = 250015050
$ws.Range("M137").Value = -83339334
$ws.Range("N137").Value = -250020150

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 23266158
$ws.Range("I32").Value = 27034158
$ws.Range("K32").Value = 27034158
$ws.Range("M32").Value = -27033871
$ws.Range("H35").Value = 1436.25
$ws.Range("I35").Value = 497.5
$ws.Range("K35").Value = 497.5
$ws.Range("M35").Value = -91.5
$ws.Range("H74").Value = 50057056
$ws.Range("I74").Value = 50057056
$ws.Range("K74").Value = 50057056
$ws.Range("M74").Value = -50056182
$ws.Range("H77").Value = 50057056
$ws.Range("I77").Value = 50057056
$ws.Range("K77").Value = 250285280
$ws.Range("M77").Value = -250280912
$ws.Range("H122").Value = 3017.4443
$ws.Range("I122").Value = 2033.6364
$ws.Range("K122").Value = 6100.9092
$ws.Range("M122").Value = -3650.9092
$ws.Range("H133").Value = 88417.8
$ws.Range("J133").Value = 96255.75
$ws.Range("L133").Value = 96255.75
$ws.Range("N133").Value = -101315.75
$ws.Range("H134").Value = 444998.25
$ws.Range("J134").Value = 444998.25
$ws.Range("L134").Value = 444998.25
$ws.Range("N134").Value = -455138.25
$ws.Range("H139").Value = 60577.57
$ws.Range("J139").Value = 60577.57
$ws.Range("L139").Value = 60577.57
$ws.Range("N139").Value = -70857.57

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 243.14285
$ws.Range("I22").Value = 164.25
$ws.Range("J22").Value = 348.33334
$ws.Range("K22").Value = 164.25
$ws.Range("L22").Value = 348.33334
$ws.Range("M22").Value = 8.75
$ws.Range("N22").Value = -694.33334
$ws.Range("H137").Value = 193593.33
$ws.Range("J137").Value = 193593.33
$ws.Range("L137").Value = 193593.33
$ws.Range("N137").Value = -203793.33

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 30309000
$ws.Range("I31").Value = 5099.9414
$ws.Range("J31").Value = 62506896
$ws.Range("K31").Value = 5099.9414
$ws.Range("L31").Value = 62506896
$ws.Range("M31").Value = -4804.9414
$ws.Range("N31").Value = -62507486
$ws.Range("H34").Value = 30309000
$ws.Range("I34").Value = 5099.9414
$ws.Range("J34").Value = 62506896
$ws.Range("K34").Value = 5099.9414
$ws.Range("L34").Value = 62506896
$ws.Range("M34").Value = -4897.9414
$ws.Range("N34").Value = -62507300
$ws.Range("H51").Value = 49995
$ws.Range("I51").Value = 49995
$ws.Range("K51").Value = 49995
$ws.Range("M51").Value = -49259
$ws.Range("H52").Value = 149984.75
$ws.Range("J52").Value = 149984.75
$ws.Range("L52").Value = 149984.75
$ws.Range("N52").Value = -150572.75
$ws.Range("H60").Value = 47838.8
$ws.Range("I60").Value = 14999.5
$ws.Range("K60").Value = 14999.5
$ws.Range("M60").Value = -14488.5
$ws.Range("H61").Value = 49995
$ws.Range("I61").Value = 49995
$ws.Range("K61").Value = 49995
$ws.Range("M61").Value = -49647
$ws.Range("H105").Value = 10917.417
$ws.Range("I105").Value = 2818
$ws.Range("K105").Value = 2818
$ws.Range("M105").Value = -1071
$ws.Range("H107").Value = 2174.4
$ws.Range("I107").Value = 1317.7273
$ws.Range("J107").Value = 3221.4443
$ws.Range("K107").Value = 1317.7273
$ws.Range("L107").Value = 3221.4443
$ws.Range("M107").Value = 602.2727
$ws.Range("N107").Value = -7061.4443
$ws.Range("H135").Value = 75419.664
$ws.Range("J135").Value = 75419.664
$ws.Range("L135").Value = 75419.664
$ws.Range("N135").Value = -85559.664
$ws.Range("H140").Value = 62372.5
$ws.Range("J140").Value = 62372.5
$ws.Range("L140").Value = 62372.5
$ws.Range("N140").Value = -72732.5
$ws.Range("H141").Value = 268885.84
$ws.Range("J141").Value = 276330.16
$ws.Range("L141").Value = 276330.16
$ws.Range("N141").Value = -286690.16

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H26").Value = 821
$ws.Range("I26").Value = 52.5
$ws.Range("K26").Value = 157.5
$ws.Range("M26").Value = 130.5
$ws.Range("H121").Value = 1030.8
$ws.Range("I121").Value = 1030.8
$ws.Range("K121").Value = 3092.4
$ws.Range("M121").Value = -1782.4
$ws.Range("H131").Value = 1749.8334
$ws.Range("J131").Value = 1866.4166
$ws.Range("L131").Value = 5599.2498
$ws.Range("N131").Value = -15679.2498
$ws.Range("H134").Value = 4710.0586
$ws.Range("I134").Value = 1540.9286
$ws.Range("J134").Value = 19499.334
$ws.Range("K134").Value = 4622.7858
$ws.Range("L134").Value = 58498.00199999999
$ws.Range("M134").Value = 447.2142000000003
$ws.Range("N134").Value = -68638.002
$ws.Range("H140").Value = 1513.0294
$ws.Range("I140").Value = 1138.3462
$ws.Range("K140").Value = 3415.0386
$ws.Range("M140").Value = 1764.9614

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 8263.917
$ws.Range("I113").Value = 7994.5
$ws.Range("K113").Value = 7994.5
$ws.Range("M113").Value = -5824.5
$ws.Range("H122").Value = 2820.6667
$ws.Range("I122").Value = 2141.6667
$ws.Range("J122").Value = 3499.6667
$ws.Range("K122").Value = 6425.000100000001
$ws.Range("L122").Value = 10499.0001
$ws.Range("M122").Value = -3975.000100000001
$ws.Range("N122").Value = -15399.0001
$ws.Range("H126").Value = 5560593
$ws.Range("I126").Value = 2946260.8
$ws.Range("K126").Value = 8838782.399999999
$ws.Range("M126").Value = -8836312.399999999
$ws.Range("H135").Value = 86857.164
$ws.Range("J135").Value = 86857.164
$ws.Range("L135").Value = 86857.164
$ws.Range("N135").Value = -96997.164

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H45").Value = 35000
$ws.Range("J45").Value = 35000
$ws.Range("L45").Value = 35000
$ws.Range("N45").Value = -35814
$ws.Range("H48").Value = 37495
$ws.Range("J48").Value = 37495
$ws.Range("L48").Value = 37495
$ws.Range("N48").Value = -38817
$ws.Range("H68").Value = 4829.4
$ws.Range("I68").Value = 3716.5
$ws.Range("J68").Value = 6498.75
$ws.Range("K68").Value = 3716.5
$ws.Range("L68").Value = 6498.75
$ws.Range("M68").Value = -2967.5
$ws.Range("N68").Value = -7996.75
$ws.Range("H71").Value = 4829.4
$ws.Range("I71").Value = 3716.5
$ws.Range("J71").Value = 6498.75
$ws.Range("K71").Value = 18582.5
$ws.Range("L71").Value = 32493.75
$ws.Range("M71").Value = -14838.5
$ws.Range("N71").Value = -39981.75
$ws.Range("H132").Value = 250005020
$ws.Range("I132").Value = 6266.3335
$ws.Range("K132").Value = 18799.0005
$ws.Range("M132").Value = -16269.0005

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 34996.668
$ws.Range("I41").Value = 34995
$ws.Range("K41").Value = 34995
$ws.Range("M41").Value = -34605
$ws.Range("H135").Value = 16749804
$ws.Range("J135").Value = 16749804
$ws.Range("L135").Value = 16749804
$ws.Range("N135").Value = -16759944
$ws.Range("H140").Value = 100000
$ws.Range("J140").Value = 100000
$ws.Range("L140").Value = 100000
$ws.Range("N140").Value = -110360
$ws.Range("H141").Value = 120000
$ws.Range("J141").Value = 120000
$ws.Range("L141").Value = 120000
$ws.Range("N141").Value = -130360

